# Do not authorize addTag for FILE object type.
# Add two new "addTag" rules:
#   1) Complaint - Anybody can add tag  (inserted before the "Case File - Assignee
#      Read Access" row, i.e. right after the last existing COMPLAINT row)
#   2) Case File - anyone can add tag   (inserted right after "Case File - anyone
#      can subscribe", i.e. right before the TASK section starts)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert the COMPLAINT "add tag" rule at row 30 -------------------------
# This pushes the existing CASE_FILE / TASK / FOLDER rows down by one row.
$ws.Rows.Item(30).Insert()

$ws.Cells.Item(30, 2).Value = "Complaint – Anybody can add tag"
$ws.Cells.Item(30, 3).Value = "COMPLAINT"
$ws.Cells.Item(30, 7).Value = "grant addTag to *"
$ws.Rows.Item(30).RowHeight = 30

# --- Insert the CASE_FILE "add tag" rule right after row 45 ----------------
# (row 45 is now "Case File - anyone can subscribe" after the shift above)
$ws.Rows.Item(46).Insert()

$ws.Cells.Item(46, 2).Value = "Case File – anyone can add tag"
$ws.Cells.Item(46, 3).Value = "CASE_FILE"
$ws.Cells.Item(46, 7).Value = "grant addTag to *"
$ws.Rows.Item(46).RowHeight = 30

# --- Update the used-range dimension / selection to match the new layout ---
$ws.Range("B62").Select()
